# "Updates to prop never treated and coverage times"
#
# Core data change: on the "Platform Coverage" sheet, row 2 (the MDA /
# "All ages 5-15" treatment coverage row) only had a coverage value every
# other year (H, J, L, N, ... every even-offset column). This fills in the
# previously-blank in-between year columns (I, K, M, O, Q, S, U, W, Y, AA,
# AC) with the same 0.6 coverage value so every year from 2018-2040 has an
# explicit "proportion never treated"/coverage figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

$colsToFill = @("I","K","M","O","Q","S","U","W","Y","AA","AC")
foreach ($col in $colsToFill) {
    $ws.Range($col + "2").Value = 0.6
}

# Leave the sheet's view positioned/selected the way it ended up after the
# edit (scrolled right, with the next empty cell after the coverage table
# selected).
[void]$ws.Activate()
[void]$ws.Range("AE2").Select()
